# Apply the commit's changes:
#  1) Rename the two worksheets.
#  2) Swap the "Prime – Advantage" / "Lifetime" pair in the health-insurance list
#     (row 27 becomes "Lifetime", row 28 becomes "Prime – Advantage").
#  3) Replace "Arogya Supreme" (row 44) with "Flexi Health" and
#     "Flexi Health" (row 45) with "Health QuBE Basic".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("sheet1")
$ws1.Name = "Student Insurance Plan"

$ws2 = $wb.Worksheets.Item("sheet2")
$ws2.Name = "All Health Insurance List"

$ws2.Cells.Item(27, 1).Value = "Lifetime"
$ws2.Cells.Item(28, 1).Value = "Prime – Advantage"

$ws2.Cells.Item(44, 1).Value = "Flexi Health"
$ws2.Cells.Item(45, 1).Value = "Health QuBE Basic"
